# Update cryptos list: Sun May  7 07:46:38 UTC 2023 (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.989.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4591"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3826"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07717"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9806"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.10"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.897.58"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.675"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.938"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07019"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.00"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009477"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.968.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.321"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.089"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.671"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.55"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09278"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8669"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.075"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.252"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.026"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05747"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.157"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5518"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.432"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.878"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.326"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5189"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.31"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06850"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.064"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000002594"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.07"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.783"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2862"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.12%  "
